$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Capture the existing "Notes" column (currently column P / index 16) ---
# including its header, before we touch anything, so we can relocate it to the
# new last column once the two new columns are appended.
$notesVals = @()
$notesIsString = @()
for ($r = 1; $r -le 23; $r++) {
    $cell = $ws.Cells.Item($r, 16)
    $notesVals += ,$cell.Value()
}

# --- Add the two new table columns (RMSLE.TissuePC, N.TissuePC) ---
# The table engine always appends new ListColumns at the end, so they land at
# columns Q (17) and R (18) for now; we fix up final left-to-right order below.
$colTissuePC = $lo.ListColumns.Add()
$colNTissuePC = $lo.ListColumns.Add()

# --- Relocate "Notes" data from P to the real last column (R / 18) ---
for ($r = 1; $r -le 23; $r++) {
    $ws.Cells.Item($r, 18).Value = $notesVals[$r - 1]
}

# --- Header row: P/Q become the new metrics, R becomes "Notes" again ---
$ws.Cells.Item(1, 16).Value = "RMSLE.TissuePC"
$ws.Cells.Item(1, 17).Value = "N.TissuePC"
$ws.Cells.Item(1, 18).Value = "Notes"

# --- New data values for RMSLE.TissuePC (P) / N.TissuePC (Q), rows 2-23 ---
$ws.Cells.Item(4, 16).Value = 0.4612
$ws.Cells.Item(4, 17).Value = 12

$ws.Cells.Item(5, 16).Value = 0.4612
$ws.Cells.Item(5, 17).Value = 12

$ws.Cells.Item(6, 16).Value = 0.5563
$ws.Cells.Item(6, 17).Value = 412

$ws.Cells.Item(7, 16).Value = 0.5925
$ws.Cells.Item(7, 17).Value = 964

$ws.Cells.Item(8, 16).Value = 0.5926
$ws.Cells.Item(8, 17).Value = 964

$ws.Cells.Item(9, 16).Value = 0.5925
$ws.Cells.Item(9, 17).Value = 964

$ws.Cells.Item(10, 16).Value = 0.6136
$ws.Cells.Item(10, 17).Value = 953

$ws.Cells.Item(11, 16).Value = 0.6136
$ws.Cells.Item(11, 17).Value = 953

$ws.Cells.Item(12, 16).Value = 0.6136
$ws.Cells.Item(12, 17).Value = 953

$ws.Cells.Item(13, 16).Value = 0.6115
$ws.Cells.Item(13, 17).Value = 964

$ws.Cells.Item(14, 16).Value = 0.6115
$ws.Cells.Item(14, 17).Value = 964

$ws.Cells.Item(15, 16).Value = 0.6098
$ws.Cells.Item(15, 17).Value = 858

$ws.Cells.Item(16, 16).Value = 0.7611
$ws.Cells.Item(16, 17).Value = 858

$ws.Cells.Item(17, 16).Value = 0.7611
$ws.Cells.Item(17, 17).Value = 858

$ws.Cells.Item(18, 16).Value = 0.7854
$ws.Cells.Item(18, 17).Value = 851

$ws.Cells.Item(19, 16).Value = 0.7866
$ws.Cells.Item(19, 17).Value = 840

$ws.Cells.Item(20, 16).Value = 0.5995
$ws.Cells.Item(20, 17).Value = 863

$ws.Cells.Item(21, 16).Value = 0.6428
$ws.Cells.Item(21, 17).Value = 863

$ws.Cells.Item(22, 16).Value = 0.643
$ws.Cells.Item(22, 17).Value = 863

$ws.Cells.Item(23, 16).Value = 0.63
$ws.Cells.Item(23, 17).Value = 863

# --- Column widths: P/Q pick up the same narrow width as the other metric
# columns (K:O), and R (now "Notes") keeps the wide "notes" column width. ---
$ws.Columns.Item(16).ColumnWidth = 5.83
$ws.Columns.Item(17).ColumnWidth = 5.83
$ws.Columns.Item(18).ColumnWidth = 38.67

# --- sheetView: reflect the selected cell after the edit (R4) ---
$ws.Range("R4").Select()
